# Applies the cryptos-list price/volume refresh described in the commit
# "Updated cryptos list on Fri Sep 15 07:40:14 UTC 2023 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.648.56"
$ws.Range("E2").Value = "  +1.28%  "
$ws.Range("D3").Value = "1.632.61"
$ws.Range("E3").Value = "  +0.81%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.38%  "
$ws.Range("E6").Value = "  +2.82%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("E8").Value = "  +2.00%  "
$ws.Range("E9").Value = "  +0.94%  "
$ws.Range("E10").Value = "  +1.56%  "
$ws.Range("E11").Value = "  +3.52%  "
$ws.Range("D12").Value = "1.860.52"
$ws.Range("E12").Value = "  +0.87%  "
$ws.Range("D13").Value = "1.633.57"
$ws.Range("E13").Value = "  +0.85%  "
$ws.Range("E14").Value = "  +1.56%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.523"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.90%  "
$ws.Range("D16").Value = "26.642.58"
$ws.Range("E16").Value = "  +1.27%  "
$ws.Range("E17").Value = "  +1.17%  "
$ws.Range("E18").Value = "  +1.73%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "218.15"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +8.04%  "
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("E22").Value = "  +2.40%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.37"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.37%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.69%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "148.05"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.51%  "
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("E27").Value = "  +1.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.86"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.51%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.48"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.97%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0504"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.79%  "
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("E32").Value = "  +3.57%  "
$ws.Range("E33").Value = "  +1.77%  "
$ws.Range("E34").Value = "  +0.53%  "
$ws.Range("E35").Value = "  +0.23%  "
$ws.Range("D36").Value = "1.209.60"
$ws.Range("E36").Value = "  +2.69%  "
$ws.Range("E37").Value = "  +5.42%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.808"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.10%  "
$ws.Range("E39").Value = "  +0.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.500"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.84%  "
$ws.Range("E41").Value = "  -1.64%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.40"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.75%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.789"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.60%  "
$ws.Range("D44").Value = "1.773.91"
$ws.Range("E44").Value = "  +0.99%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.91"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.41%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.55"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.51%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "54.71"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.62%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0513"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.60"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.16%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.409"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.09%  "
$ws.Range("E51").Value = "  +0.27%  "
